$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Lama1"
$ws.Cells.Item(2, 3).Value = "Itgb1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.014112
$ws.Cells.Item(2, 8).Value = 0.042336
$ws.Cells.Item(2, 9).Value = 0.1773673913134555
$ws.Cells.Item(2, 10).Value = 0.1773673913134555
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 176.8550973333333
$ws.Cells.Item(2, 14).Value = 530.565292
$ws.Cells.Item(2, 15).Value = 0.2669710696905332
$ws.Cells.Item(2, 16).Value = 0.2669710696905332
$ws.Cells.Item(2, 17).Value = 2.495779133568
$ws.Cells.Item(2, 18).Value = 22.462012202112
$ws.Cells.Item(2, 19).Value = 0.0473519621871726
$ws.Cells.Item(2, 20).Value = 0.04735196218717259

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Lama1"
$ws.Cells.Item(3, 3).Value = "Itgb1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.014112
$ws.Cells.Item(3, 8).Value = 0.042336
$ws.Cells.Item(3, 9).Value = 0.1773673913134555
$ws.Cells.Item(3, 10).Value = 0.1773673913134555
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 145.6413626666667
$ws.Cells.Item(3, 14).Value = 436.924088
$ws.Cells.Item(3, 15).Value = 0.2198524722701247
$ws.Cells.Item(3, 16).Value = 0.2198524722701247
$ws.Cells.Item(3, 17).Value = 2.055290909952
$ws.Cells.Item(3, 18).Value = 18.497618189568
$ws.Cells.Item(3, 19).Value = 0.03899465948036584
$ws.Cells.Item(3, 20).Value = 0.03899465948036583

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Lama1"
$ws.Cells.Item(4, 3).Value = "Itgb1"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.014112
$ws.Cells.Item(4, 8).Value = 0.042336
$ws.Cells.Item(4, 9).Value = 0.1773673913134555
$ws.Cells.Item(4, 10).Value = 0.1773673913134555
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 84.02511333333334
$ws.Cells.Item(4, 14).Value = 252.07534
$ws.Cells.Item(4, 15).Value = 0.1268398520919549
$ws.Cells.Item(4, 16).Value = 0.1268398520919549
$ws.Cells.Item(4, 17).Value = 1.18576239936
$ws.Cells.Item(4, 18).Value = 10.67186159424
$ws.Cells.Item(4, 19).Value = 0.02249725368013458
$ws.Cells.Item(4, 20).Value = 0.02249725368013457

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Lama1"
$ws.Cells.Item(5, 3).Value = "Itgb1"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.014112
$ws.Cells.Item(5, 8).Value = 0.042336
$ws.Cells.Item(5, 9).Value = 0.1773673913134555
$ws.Cells.Item(5, 10).Value = 0.1773673913134555
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 84.92877566666668
$ws.Cells.Item(5, 14).Value = 254.786327
$ws.Cells.Item(5, 15).Value = 0.1282039727953256
$ws.Cells.Item(5, 16).Value = 0.1282039727953256
$ws.Cells.Item(5, 17).Value = 1.198514882208
$ws.Cells.Item(5, 18).Value = 10.786633939872
$ws.Cells.Item(5, 19).Value = 0.02273920421072812
$ws.Cells.Item(5, 20).Value = 0.02273920421072812

# Row 6
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Lama1"
$ws.Cells.Item(6, 3).Value = "Itgb1"
$ws.Cells.Item(6, 4).Value = "Neutrophils"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.014112
$ws.Cells.Item(6, 8).Value = 0.042336
$ws.Cells.Item(6, 9).Value = 0.1773673913134555
$ws.Cells.Item(6, 10).Value = 0.1773673913134555
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 63.97102366666667
$ws.Cells.Item(6, 14).Value = 191.913071
$ws.Cells.Item(6, 15).Value = 0.09656726254996952
$ws.Cells.Item(6, 16).Value = 0.09656726254996952
$ws.Cells.Item(6, 17).Value = 0.902759085984
$ws.Cells.Item(6, 18).Value = 8.124831773856
$ws.Cells.Item(6, 19).Value = 0.01712788344476964
$ws.Cells.Item(6, 20).Value = 0.01712788344476964

# Row 7
$ws.Cells.Item(7, 1).Value = "ECs"
$ws.Cells.Item(7, 2).Value = "Lama1"
$ws.Cells.Item(7, 3).Value = "Itgb1"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.014112
$ws.Cells.Item(7, 8).Value = 0.042336
$ws.Cells.Item(7, 9).Value = 0.1773673913134555
$ws.Cells.Item(7, 10).Value = 0.1773673913134555
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 107.0290476666667
$ws.Cells.Item(7, 14).Value = 321.087143
$ws.Cells.Item(7, 15).Value = 0.1615653706020921
$ws.Cells.Item(7, 16).Value = 0.1615653706020921
$ws.Cells.Item(7, 17).Value = 1.510393920672
$ws.Cells.Item(7, 18).Value = 13.593545286048
$ws.Cells.Item(7, 19).Value = 0.02865642831028473
$ws.Cells.Item(7, 20).Value = 0.02865642831028472

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Lama1"
$ws.Cells.Item(8, 3).Value = "Itgb1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.059263
$ws.Cells.Item(8, 8).Value = 0.177789
$ws.Cells.Item(8, 9).Value = 0.7448500362393221
$ws.Cells.Item(8, 10).Value = 0.7448500362393219
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 176.8550973333333
$ws.Cells.Item(8, 14).Value = 530.565292
$ws.Cells.Item(8, 15).Value = 0.2669710696905332
$ws.Cells.Item(8, 16).Value = 0.2669710696905332
$ws.Cells.Item(8, 17).Value = 10.48096363326533
$ws.Cells.Item(8, 18).Value = 94.32867269938801
$ws.Cells.Item(8, 19).Value = 0.1988534109338442
$ws.Cells.Item(8, 20).Value = 0.1988534109338442

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Lama1"
$ws.Cells.Item(9, 3).Value = "Itgb1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.059263
$ws.Cells.Item(9, 8).Value = 0.177789
$ws.Cells.Item(9, 9).Value = 0.7448500362393221
$ws.Cells.Item(9, 10).Value = 0.7448500362393219
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 145.6413626666667
$ws.Cells.Item(9, 14).Value = 436.924088
$ws.Cells.Item(9, 15).Value = 0.2198524722701247
$ws.Cells.Item(9, 16).Value = 0.2198524722701247
$ws.Cells.Item(9, 17).Value = 8.631144075714666
$ws.Cells.Item(9, 18).Value = 77.680296681432
$ws.Cells.Item(9, 19).Value = 0.163757121937707
$ws.Cells.Item(9, 20).Value = 0.1637571219377069

# Row 10
$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Lama1"
$ws.Cells.Item(10, 3).Value = "Itgb1"
$ws.Cells.Item(10, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.059263
$ws.Cells.Item(10, 8).Value = 0.177789
$ws.Cells.Item(10, 9).Value = 0.7448500362393221
$ws.Cells.Item(10, 10).Value = 0.7448500362393219
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 84.02511333333334
$ws.Cells.Item(10, 14).Value = 252.07534
$ws.Cells.Item(10, 15).Value = 0.1268398520919549
$ws.Cells.Item(10, 16).Value = 0.1268398520919549
$ws.Cells.Item(10, 17).Value = 4.979580291473334
$ws.Cells.Item(10, 18).Value = 44.81622262326
$ws.Cells.Item(10, 19).Value = 0.09447666842728283
$ws.Cells.Item(10, 20).Value = 0.09447666842728281

# Row 11
$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Lama1"
$ws.Cells.Item(11, 3).Value = "Itgb1"
$ws.Cells.Item(11, 4).Value = "MuSCs"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.059263
$ws.Cells.Item(11, 8).Value = 0.177789
$ws.Cells.Item(11, 9).Value = 0.7448500362393221
$ws.Cells.Item(11, 10).Value = 0.7448500362393219
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 84.92877566666668
$ws.Cells.Item(11, 14).Value = 254.786327
$ws.Cells.Item(11, 15).Value = 0.1282039727953256
$ws.Cells.Item(11, 16).Value = 0.1282039727953256
$ws.Cells.Item(11, 17).Value = 5.033134032333668
$ws.Cells.Item(11, 18).Value = 45.29820629100301
$ws.Cells.Item(11, 19).Value = 0.09549273378262335
$ws.Cells.Item(11, 20).Value = 0.09549273378262332

# Row 12
$ws.Cells.Item(12, 1).Value = "FAPs"
$ws.Cells.Item(12, 2).Value = "Lama1"
$ws.Cells.Item(12, 3).Value = "Itgb1"
$ws.Cells.Item(12, 4).Value = "Neutrophils"
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.059263
$ws.Cells.Item(12, 8).Value = 0.177789
$ws.Cells.Item(12, 9).Value = 0.7448500362393221
$ws.Cells.Item(12, 10).Value = 0.7448500362393219
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 63.97102366666667
$ws.Cells.Item(12, 14).Value = 191.913071
$ws.Cells.Item(12, 15).Value = 0.09656726254996952
$ws.Cells.Item(12, 16).Value = 0.09656726254996952
$ws.Cells.Item(12, 17).Value = 3.791114775557667
$ws.Cells.Item(12, 18).Value = 34.120032980019
$ws.Cells.Item(12, 19).Value = 0.07192812900987693
$ws.Cells.Item(12, 20).Value = 0.0719281290098769

# Row 13
$ws.Cells.Item(13, 1).Value = "FAPs"
$ws.Cells.Item(13, 2).Value = "Lama1"
$ws.Cells.Item(13, 3).Value = "Itgb1"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.059263
$ws.Cells.Item(13, 8).Value = 0.177789
$ws.Cells.Item(13, 9).Value = 0.7448500362393221
$ws.Cells.Item(13, 10).Value = 0.7448500362393219
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 107.0290476666667
$ws.Cells.Item(13, 14).Value = 321.087143
$ws.Cells.Item(13, 15).Value = 0.1615653706020921
$ws.Cells.Item(13, 16).Value = 0.1615653706020921
$ws.Cells.Item(13, 17).Value = 6.342862451869667
$ws.Cells.Item(13, 18).Value = 57.08576206682701
$ws.Cells.Item(13, 19).Value = 0.1203419721479878
$ws.Cells.Item(13, 20).Value = 0.1203419721479878

# Row 14
$ws.Cells.Item(14, 1).Value = "MuSCs"
$ws.Cells.Item(14, 2).Value = "Lama1"
$ws.Cells.Item(14, 3).Value = "Itgb1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.006188666666666666
$ws.Cells.Item(14, 8).Value = 0.018566
$ws.Cells.Item(14, 9).Value = 0.07778257244722256
$ws.Cells.Item(14, 10).Value = 0.07778257244722254
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 176.8550973333333
$ws.Cells.Item(14, 14).Value = 530.565292
$ws.Cells.Item(14, 15).Value = 0.2669710696905332
$ws.Cells.Item(14, 16).Value = 0.2669710696905332
$ws.Cells.Item(14, 17).Value = 1.094497245696889
$ws.Cells.Item(14, 18).Value = 9.850475211272
$ws.Cells.Item(14, 19).Value = 0.0207656965695164
$ws.Cells.Item(14, 20).Value = 0.0207656965695164

# Row 15
$ws.Cells.Item(15, 1).Value = "MuSCs"
$ws.Cells.Item(15, 2).Value = "Lama1"
$ws.Cells.Item(15, 3).Value = "Itgb1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.006188666666666666
$ws.Cells.Item(15, 8).Value = 0.018566
$ws.Cells.Item(15, 9).Value = 0.07778257244722256
$ws.Cells.Item(15, 10).Value = 0.07778257244722254
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 145.6413626666667
$ws.Cells.Item(15, 14).Value = 436.924088
$ws.Cells.Item(15, 15).Value = 0.2198524722701247
$ws.Cells.Item(15, 16).Value = 0.2198524722701247
$ws.Cells.Item(15, 17).Value = 0.901325846423111
$ws.Cells.Item(15, 18).Value = 8.111932617808
$ws.Cells.Item(15, 19).Value = 0.01710069085205197
$ws.Cells.Item(15, 20).Value = 0.01710069085205197

# Row 16
$ws.Cells.Item(16, 1).Value = "MuSCs"
$ws.Cells.Item(16, 2).Value = "Lama1"
$ws.Cells.Item(16, 3).Value = "Itgb1"
$ws.Cells.Item(16, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.006188666666666666
$ws.Cells.Item(16, 8).Value = 0.018566
$ws.Cells.Item(16, 9).Value = 0.07778257244722256
$ws.Cells.Item(16, 10).Value = 0.07778257244722254
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 84.02511333333334
$ws.Cells.Item(16, 14).Value = 252.07534
$ws.Cells.Item(16, 15).Value = 0.1268398520919549
$ws.Cells.Item(16, 16).Value = 0.1268398520919549
$ws.Cells.Item(16, 17).Value = 0.5200034180488888
$ws.Cells.Item(16, 18).Value = 4.68003076244
$ws.Cells.Item(16, 19).Value = 0.009865929984537473
$ws.Cells.Item(16, 20).Value = 0.009865929984537471

# Row 17
$ws.Cells.Item(17, 1).Value = "MuSCs"
$ws.Cells.Item(17, 2).Value = "Lama1"
$ws.Cells.Item(17, 3).Value = "Itgb1"
$ws.Cells.Item(17, 4).Value = "MuSCs"
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.006188666666666666
$ws.Cells.Item(17, 8).Value = 0.018566
$ws.Cells.Item(17, 9).Value = 0.07778257244722256
$ws.Cells.Item(17, 10).Value = 0.07778257244722254
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 84.92877566666668
$ws.Cells.Item(17, 14).Value = 254.786327
$ws.Cells.Item(17, 15).Value = 0.1282039727953256
$ws.Cells.Item(17, 16).Value = 0.1282039727953256
$ws.Cells.Item(17, 17).Value = 0.5255958830091112
$ws.Cells.Item(17, 18).Value = 4.730362947082
$ws.Cells.Item(17, 19).Value = 0.009972034801974164
$ws.Cells.Item(17, 20).Value = 0.009972034801974162

# Row 18
$ws.Cells.Item(18, 1).Value = "MuSCs"
$ws.Cells.Item(18, 2).Value = "Lama1"
$ws.Cells.Item(18, 3).Value = "Itgb1"
$ws.Cells.Item(18, 4).Value = "Neutrophils"
$ws.Cells.Item(18, 5).Value = 1
$ws.Cells.Item(18, 6).Value = 0.3333333333333333
$ws.Cells.Item(18, 7).Value = 0.006188666666666666
$ws.Cells.Item(18, 8).Value = 0.018566
$ws.Cells.Item(18, 9).Value = 0.07778257244722256
$ws.Cells.Item(18, 10).Value = 0.07778257244722254
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 63.97102366666667
$ws.Cells.Item(18, 14).Value = 191.913071
$ws.Cells.Item(18, 15).Value = 0.09656726254996952
$ws.Cells.Item(18, 16).Value = 0.09656726254996952
$ws.Cells.Item(18, 17).Value = 0.3958953417984444
$ws.Cells.Item(18, 18).Value = 3.563058076186
$ws.Cells.Item(18, 19).Value = 0.007511250095322965
$ws.Cells.Item(18, 20).Value = 0.007511250095322964

# Row 19
$ws.Cells.Item(19, 1).Value = "MuSCs"
$ws.Cells.Item(19, 2).Value = "Lama1"
$ws.Cells.Item(19, 3).Value = "Itgb1"
$ws.Cells.Item(19, 4).Value = "Resolving-Mac"
$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 6).Value = 0.3333333333333333
$ws.Cells.Item(19, 7).Value = 0.006188666666666666
$ws.Cells.Item(19, 8).Value = 0.018566
$ws.Cells.Item(19, 9).Value = 0.07778257244722256
$ws.Cells.Item(19, 10).Value = 0.07778257244722254
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 107.0290476666667
$ws.Cells.Item(19, 14).Value = 321.087143
$ws.Cells.Item(19, 15).Value = 0.1615653706020921
$ws.Cells.Item(19, 16).Value = 0.1615653706020921
$ws.Cells.Item(19, 17).Value = 0.6623670996597777
$ws.Cells.Item(19, 18).Value = 5.961303896938
$ws.Cells.Item(19, 19).Value = 0.01256697014381959
$ws.Cells.Item(19, 20).Value = 0.01256697014381959
